$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162, shifting the existing rows 162..228 down to 163..229.
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly price-report entry
# (Feria Lagunitas de Puerto Montt - Acelga).
$ws.Range("A162").Value = 4
$ws.Range("B162").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C162").Value = "Los Lagos"
$ws.Range("D162").Value = 44845
$ws.Range("E162").Value = 10
$ws.Range("F162").Value = 100112009
$ws.Range("G162").Value = "Acelga"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 180
$ws.Range("K162").Value = 3000
$ws.Range("L162").Value = 3500
$ws.Range("M162").Value = 3222
$ws.Range("N162").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O162").Value = "Región del Maule"
$ws.Range("P162").Value = 806
$ws.Range("Q162").Value = 4
$ws.Range("R162").Value = "Hortaliza"
